$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10 (ALC)
$ws.Range("H10").Value = 3111.5557
$ws.Range("J10").Value = 3063.125
$ws.Range("L10").Value = 3063.125
$ws.Range("N10").Value = -3649.125

# Row 92 (ALC)
$ws.Range("H92").Value = 1268.0435
$ws.Range("I92").Value = 1240.75
$ws.Range("J92").Value = 1450
$ws.Range("K92").Value = 1240.75
$ws.Range("L92").Value = 1450
$ws.Range("M92").Value = 7.25
$ws.Range("N92").Value = -3946

# Row 137 (ALC)
$ws.Range("H137").Value = 821.7143
$ws.Range("I137").Value = 747.6111
$ws.Range("J137").Value = 1266.3334
$ws.Range("K137").Value = 2242.8333
$ws.Range("L137").Value = 3799.0002
$ws.Range("M137").Value = 307.1667000000002
$ws.Range("N137").Value = -8899.0002

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (ARM)
$ws.Range("H61").Value = 1116.6666
$ws.Range("I61").Value = 940
$ws.Range("K61").Value = 940
$ws.Range("M61").Value = -728

# Row 62 (ARM)
$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9376

# Row 65 (ARM)
$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 30000
$ws.Range("M65").Value = -26880

# Row 122 (ARM)
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# Row 136 (ARM)
$ws.Range("H136").Value = 1116.6666
$ws.Range("I136").Value = 940
$ws.Range("K136").Value = 2820
$ws.Range("M136").Value = -270

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (CRP)
$ws.Range("H4").Value = 1490
$ws.Range("J4").Value = 1490
$ws.Range("L4").Value = 1490
$ws.Range("N4").Value = -1714

# Row 31 (CRP)
$ws.Range("H31").Value = 2176.3447
$ws.Range("I31").Value = 2078.88
$ws.Range("J31").Value = 2785.5
$ws.Range("K31").Value = 2078.88
$ws.Range("L31").Value = 2785.5
$ws.Range("M31").Value = -1783.88
$ws.Range("N31").Value = -3375.5

# Row 34 (CRP)
$ws.Range("H34").Value = 2176.3447
$ws.Range("I34").Value = 2078.88
$ws.Range("J34").Value = 2785.5
$ws.Range("K34").Value = 2078.88
$ws.Range("L34").Value = 2785.5
$ws.Range("M34").Value = -1876.88
$ws.Range("N34").Value = -3189.5

# Row 58 (CRP)
$ws.Range("H58").Value = 3823.8823
$ws.Range("I58").Value = 846.3
$ws.Range("K58").Value = 846.3
$ws.Range("M58").Value = -643.3

# Row 62 (CRP)
$ws.Range("H62").Value = 11160

# Row 65 (CRP)
$ws.Range("H65").Value = 11160

# Row 118 (CRP)
$ws.Range("H118").Value = 54900
$ws.Range("J118").Value = 54900
$ws.Range("L118").Value = 54900
$ws.Range("N118").Value = -58214

# Row 136 (CRP)
$ws.Range("H136").Value = 3823.8823
$ws.Range("I136").Value = 846.3
$ws.Range("K136").Value = 2538.9
$ws.Range("M136").Value = 11.10000000000036

$ws = $wb.Worksheets.Item("CUL")
# Row 3 (CUL)
$ws.Range("H3").Value = 6312.5
$ws.Range("I3").Value = 2666.6667
$ws.Range("K3").Value = 8000.000100000001
$ws.Range("M3").Value = -7888.000100000001

# Row 63 (CUL)
$ws.Range("H63").Value = 4766.5
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 4927.091
$ws.Range("K63").Value = 9000
$ws.Range("L63").Value = 14781.273
$ws.Range("M63").Value = -8251
$ws.Range("N63").Value = -16279.273

# Row 64 (CUL)
$ws.Range("H64").Value = 2301.182
$ws.Range("J64").Value = 2914.125
$ws.Range("L64").Value = 8742.375
$ws.Range("N64").Value = -9282.375

# Row 66 (CUL)
$ws.Range("H66").Value = 4766.5
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 4927.091
$ws.Range("K66").Value = 27000
$ws.Range("L66").Value = 44343.819
$ws.Range("M66").Value = -23256
$ws.Range("N66").Value = -51831.819

# Row 67 (CUL)
$ws.Range("H67").Value = 2301.182
$ws.Range("J67").Value = 2914.125
$ws.Range("L67").Value = 8742.375
$ws.Range("N67").Value = -10614.375

# Row 87 (CUL)
$ws.Range("H87").Value = 13583.333
$ws.Range("I87").Value = 2760
$ws.Range("J87").Value = 21314.285
$ws.Range("K87").Value = 8280
$ws.Range("L87").Value = 63942.855
$ws.Range("M87").Value = -7032
$ws.Range("N87").Value = -66438.855

# Row 90 (CUL)
$ws.Range("H90").Value = 13583.333
$ws.Range("I90").Value = 2760
$ws.Range("J90").Value = 21314.285
$ws.Range("K90").Value = 24840
$ws.Range("L90").Value = 191828.565
$ws.Range("M90").Value = -18600
$ws.Range("N90").Value = -204308.565

# Row 113 (CUL)
$ws.Range("H113").Value = 767.9091
$ws.Range("I113").Value = 534.3333
$ws.Range("J113").Value = 791.26666
$ws.Range("K113").Value = 1602.9999
$ws.Range("L113").Value = 2373.79998
$ws.Range("M113").Value = 567.0001
$ws.Range("N113").Value = -6713.79998

# Row 127 (CUL)
$ws.Range("H127").Value = 1066.6666
$ws.Range("J127").Value = 1066.6666
$ws.Range("L127").Value = 3199.9998
$ws.Range("N127").Value = -13119.9998

# Row 131 (CUL)
$ws.Range("H131").Value = 24657.455
$ws.Range("I131").Value = 59815.707
$ws.Range("J131").Value = 2520.7778
$ws.Range("K131").Value = 179447.121
$ws.Range("L131").Value = 7562.3334
$ws.Range("M131").Value = -174407.121
$ws.Range("N131").Value = -17642.3334

# Row 139 (CUL)
$ws.Range("H139").Value = 2460.6
$ws.Range("I139").Value = 1839.1538
$ws.Range("J139").Value = 6500
$ws.Range("K139").Value = 5517.4614
$ws.Range("L139").Value = 19500
$ws.Range("M139").Value = -377.4614000000001
$ws.Range("N139").Value = -29780

# Row 140 (CUL)
$ws.Range("H140").Value = 121098.586
$ws.Range("I140").Value = 203709.7
$ws.Range("J140").Value = 3082.7144
$ws.Range("K140").Value = 611129.1000000001
$ws.Range("L140").Value = 9248.143199999999
$ws.Range("M140").Value = -605949.1000000001
$ws.Range("N140").Value = -19608.1432

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM)
$ws.Range("H70").Value = 5383.6924
$ws.Range("I70").Value = 4888.8887
$ws.Range("J70").Value = 6497
$ws.Range("K70").Value = 4888.8887
$ws.Range("L70").Value = 6497
$ws.Range("M70").Value = -4618.8887
$ws.Range("N70").Value = -7037

# Row 73 (GSM)
$ws.Range("H73").Value = 5383.6924
$ws.Range("I73").Value = 4888.8887
$ws.Range("J73").Value = 6497
$ws.Range("K73").Value = 4888.8887
$ws.Range("L73").Value = 6497
$ws.Range("M73").Value = -3952.8887
$ws.Range("N73").Value = -8369

# Row 102 (GSM)
$ws.Range("H102").Value = 1740.9412
$ws.Range("I102").Value = 1276.6923
$ws.Range("J102").Value = 3249.75
$ws.Range("K102").Value = 1276.6923
$ws.Range("L102").Value = 3249.75
$ws.Range("M102").Value = 345.3077000000001
$ws.Range("N102").Value = -6493.75

# Row 122 (GSM)
$ws.Range("H122").Value = 2558.8572
$ws.Range("I122").Value = 2452.3333
$ws.Range("K122").Value = 7356.999899999999
$ws.Range("M122").Value = -4906.999899999999

# Row 132 (GSM)
$ws.Range("H132").Value = 3635.524
$ws.Range("I132").Value = 4170.4443
$ws.Range("K132").Value = 12511.3329
$ws.Range("M132").Value = -9981.332900000001

# Row 135 (GSM)
$ws.Range("H135").Value = 27099.25
$ws.Range("J135").Value = 27099.25
$ws.Range("L135").Value = 27099.25
$ws.Range("N135").Value = -37239.25

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (LTW)
$ws.Range("H2").Value = 980
$ws.Range("J2").Value = 980
$ws.Range("L2").Value = 980
$ws.Range("N2").Value = -1204

# Row 46 (LTW)
$ws.Range("H46").Value = 1253.909
$ws.Range("I46").Value = 1513.2858
$ws.Range("J46").Value = 800
$ws.Range("K46").Value = 1513.2858
$ws.Range("L46").Value = 800
$ws.Range("M46").Value = -1325.2858
$ws.Range("N46").Value = -1176

# Row 122 (LTW)
$ws.Range("H122").Value = 5051.5
$ws.Range("I122").Value = 8227.385
$ws.Range("J122").Value = 3085.476
$ws.Range("K122").Value = 24682.155
$ws.Range("L122").Value = 9256.428
$ws.Range("M122").Value = -22232.155
$ws.Range("N122").Value = -14156.428

$ws = $wb.Worksheets.Item("WVR")
# Row 136 (WVR)
$ws.Range("H136").Value = 2151.0857
$ws.Range("I136").Value = 2163.9678
$ws.Range("J136").Value = 2051.25
$ws.Range("K136").Value = 6491.903399999999
$ws.Range("L136").Value = 6153.75
$ws.Range("M136").Value = -3941.903399999999
$ws.Range("N136").Value = -11253.75
